# Add a new "2022-Q4" quarterly sheet and refresh the "总计" (summary) sheet
# to include it, shifting the existing quarterly rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet.
#    Easiest reliable way to inherit the exact look & feel (styles,
#    column layout, borders, header row) used by every other quarterly
#    sheet is to clone an existing quarterly sheet (e.g. "2022-Q3") and
#    then overwrite its data with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($template)                 # inserted immediately before the template
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The template ("2022-Q3") has 4 data rows; 2022-Q4 only needs 3, so drop row 5.
$newSheet.Range("A5:H5").Clear()

# Header row (unchanged wording/style, just confirming contents match the template)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows. Columns B, D, E, F, G hold numeric-looking figures that must stay
# TEXT (matching the rest of the workbook), so a leading apostrophe forces
# Excel to store them as literal strings instead of coercing to numbers.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'013659"
$newSheet.Range("C2").Value = "中融金融鑫选3个月持有混合A"
$newSheet.Range("D2").Value = "'1.39"
$newSheet.Range("E2").Value = "'90.62"
$newSheet.Range("F2").Value = "'4.96"
$newSheet.Range("G2").Value = "'0.0689"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'013660"
$newSheet.Range("C3").Value = "中融金融鑫选3个月持有混合C"
$newSheet.Range("D3").Value = "'0.83"
$newSheet.Range("E3").Value = "'90.62"
$newSheet.Range("F3").Value = "'4.96"
$newSheet.Range("G3").Value = "'0.0412"
$newSheet.Range("H3").Value = 8

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'516980"
$newSheet.Range("C4").Value = "华富中证证券公司先锋策略ETF"
$newSheet.Range("D4").Value = "'0.25"
$newSheet.Range("E4").Value = "'99.41"
$newSheet.Range("F4").Value = "'2.67"
$newSheet.Range("G4").Value = "'0.0067"
$newSheet.Range("H4").Value = 8

# Restore the originally-active tab (2020-Q4, the last sheet) — copying a
# sheet makes the copy the active tab, which would otherwise change the
# workbook's saved selection state.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

# ---------------------------------------------------------------------
# 2. Refresh the "总计" summary sheet: insert the 2022-Q4 totals as the
#    new second row and push every following quarter down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.12

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 0.25

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 11
$summary.Range("D4").Value = 0.81

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 14
$summary.Range("D5").Value = 0.75

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 23
$summary.Range("D6").Value = 3.79

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 14
$summary.Range("D7").Value = 1.72

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q2"
$summary.Range("C8").Value = 1
$summary.Range("D8").Value = 0.08

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2021-Q1"
$summary.Range("C9").Value = 15
$summary.Range("D9").Value = 1.58

$summary.Range("A10").Value = 8
$summary.Range("B10").Value = "2020-Q4"
$summary.Range("C10").Value = 9
$summary.Range("D10").Value = 0.33

# Row 10 is brand new (the sheet used to stop at row 9), so A10 needs to pick
# up the same "index column" style (s="2") the rows above it already use.
$summary.Range("A9").Copy()
$summary.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
